$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156; this shifts the existing rows 156:178 down to 157:179
# and extends the used range dimension to T179, matching the target diff.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new price-record data.
$ws.Cells.Item(156, 1).Value  = 7
$ws.Cells.Item(156, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(156, 3).Value  = "Ñuble"
$ws.Cells.Item(156, 4).Value  = 44522
$ws.Cells.Item(156, 5).Value  = 16
$ws.Cells.Item(156, 6).Value  = "Fruta"
$ws.Cells.Item(156, 7).Value  = 100108
$ws.Cells.Item(156, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(156, 9).Value  = 100108005
$ws.Cells.Item(156, 10).Value = "Piña"
$ws.Cells.Item(156, 11).Value = "Caramelo"
$ws.Cells.Item(156, 12).Value = "Primera"
$ws.Cells.Item(156, 13).Value = 60
$ws.Cells.Item(156, 14).Value = 18000
$ws.Cells.Item(156, 15).Value = 19000
$ws.Cells.Item(156, 16).Value = 18500
$ws.Cells.Item(156, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(156, 18).Value = "Ecuador"
$ws.Cells.Item(156, 19).Value = 1542
$ws.Cells.Item(156, 20).Value = 12
